$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.169.86"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.55%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.588.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.78%  "
$ws.Range("E6").Value = "  -1.05%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("E9").Value = "  -1.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.95"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0844"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.812.01"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.599.39"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.69%  "
$ws.Range("E14").Value = "  -1.51%  "
$ws.Range("E15").Value = "  -1.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.50"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.32%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.164.98"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0723"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.85%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "213.73"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.32%  "
$ws.Range("E20").Value = "  -1.68%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.02"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.67%  "
$ws.Range("E24").Value = "  -1.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.73"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.43%  "
$ws.Range("E28").Value = "  -1.49%  "
$ws.Range("E29").Value = "  -1.19%  "
$ws.Range("E30").Value = "  -2.63%  "
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("E32").Value = "  -2.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.416.81"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.56%  "
$ws.Range("E34").Value = "  -2.00%  "
$ws.Range("E35").Value = "  -0.75%  "
$ws.Range("E36").Value = "  -1.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.587"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.06%  "
$ws.Range("E38").Value = "  -1.81%  "
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.88"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.49%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.820"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.84%  "
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.973"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -10.01%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.763"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.40%  "
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.13"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.723.94"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.90"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "86.77"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.01%  "
$ws.Range("E48").Value = "  -0.81%  "
$ws.Range("E49").Value = "  -0.95%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0956"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.40%  "
$ws.Range("E51").Value = "  -0.20%  "
